$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, pushing existing rows 214.. down by one
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with this week's new data point
$ws.Cells.Item(214,1).Value2  = 10
$ws.Cells.Item(214,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(214,3).Value2  = "La Araucanía"
$ws.Cells.Item(214,4).Value2  = 45015
$ws.Cells.Item(214,5).Value2  = 9
$ws.Cells.Item(214,6).Value2  = "Fruta"
$ws.Cells.Item(214,7).Value2  = 100104
$ws.Cells.Item(214,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(214,9).Value2  = 100104003
$ws.Cells.Item(214,10).Value2 = "Membrillo"
$ws.Cells.Item(214,11).Value2 = "Champion"
$ws.Cells.Item(214,12).Value2 = "Primera"
$ws.Cells.Item(214,13).Value2 = 65
$ws.Cells.Item(214,14).Value2 = 14000
$ws.Cells.Item(214,15).Value2 = 14000
$ws.Cells.Item(214,16).Value2 = 14000
$ws.Cells.Item(214,17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(214,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(214,19).Value2 = 778
$ws.Cells.Item(214,20).Value2 = 18
